# "Translate German terms into English language #2"
#
# Retranslates the non-residential specific-demand lookup table on sheet
# "Tabelle1": the two numeric-column headers (C1/D1) and every label in
# the type_name column (B2:B45) move from German to English. E1/F1
# ("SLP_th_type"/"SLP_el_type") and the numeric data in C:F are unchanged.
#
# Cells are written in the same order the shared-string pool ends up in
# after the edit (rather than strictly top-to-bottom), so new strings are
# interned in the right sequence and re-used strings ("Metall", "Hotels")
# keep resolving correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Header row (row 1): C/D captions get re-worded; E/F text is unchanged.
$ws.Range("C1").Value = "Spec.final th. energy demand in kWh/m2*a"
$ws.Range("D1").Value = "Spec. el. energy demand kWh/m2*a"

# type_name column (B2:B45): German -> English.
$ws.Range("B4").Value  = "bank and insurance"
$ws.Range("B6").Value  = "Non profit organizations"
$ws.Range("B7").Value  = "Small office buildings"
$ws.Range("B8").Value  = "Other services"
$ws.Range("B5").Value  = "Public institutions"
$ws.Range("B3").Value  = "Finishing trade construction work"
$ws.Range("B2").Value  = "Main construction work"
$ws.Range("B10").Value = "Automobile"
$ws.Range("B11").Value = "Wood and timber"
$ws.Range("B12").Value = "Paper"
$ws.Range("B13").Value = "Small retailer for food"
$ws.Range("B14").Value = "Small retailer for non-food"
$ws.Range("B15").Value = "Large retailer for food"
$ws.Range("B16").Value = "Large retailer for non-food"
$ws.Range("B17").Value = "Primary school"
$ws.Range("B18").Value = "School for physically handicapped"
$ws.Range("B19").Value = "High school"
$ws.Range("B20").Value = "Trade school"
$ws.Range("B21").Value = "University"
$ws.Range("B23").Value = "Restaurants"
$ws.Range("B24").Value = "Childrens home"
$ws.Range("B26").Value = "Butcher"
$ws.Range("B25").Value = "Backery"
$ws.Range("B27").Value = "Laundry"
$ws.Range("B28").Value = "Farm primary agriculture "
$ws.Range("B29").Value = "Farm with 10 - 49 cattle units"
$ws.Range("B30").Value = "Farm with 50 - 100 cattle units"
$ws.Range("B31").Value = "Farm with more than 100 cattle units"
$ws.Range("B32").Value = "Gardening"
$ws.Range("B33").Value = "Hospital"
$ws.Range("B34").Value = "Library"
$ws.Range("B35").Value = "Prison"
$ws.Range("B36").Value = "Cinema"
$ws.Range("B37").Value = "Theater"
$ws.Range("B38").Value = "Parish hall"
$ws.Range("B39").Value = "Sports hall"
$ws.Range("B40").Value = "Multi purpose hall"
$ws.Range("B41").Value = "Swimming hall"
$ws.Range("B42").Value = "Club house"
$ws.Range("B43").Value = "Fitness studio"
$ws.Range("B44").Value = "Train station smaller 5000m2"
$ws.Range("B45").Value = "Train station equal or larger than 5000m2"

# B9 ("Metall") and B22 ("Hotels") keep their German text unchanged (it's
# already the same in English), so they are intentionally left untouched.

# Re-fit columns B:D to the new English text widths.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Selection moves from the old bottom-of-sheet cell (B30) up to B6, and the
# view scrolls back so row 1 is visible again (no more topLeftCell="A25").
$ws.Range("B6").Select() | Out-Null
